$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update temperature values for the "Worm 3 / LGF" rows (data correction) ---
$ws.Range("C17").Value = 22.8
$ws.Range("C18").Value = 17.6
$ws.Range("C19").Value = 13.1

# --- Sort the whole data range ascending by the Temp (deg C) column ---
$dataRange = $ws.Range("A1:I19")
$sortKey = $ws.Range("C1:C19")
$tieKey = $ws.Range("B1:B19")
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($sortKey, [Microsoft.Office.Interop.Excel.XlSortOn]::xlSortOnValues, [Microsoft.Office.Interop.Excel.XlSortOrder]::xlAscending, [Microsoft.Office.Interop.Excel.XlSortDataOption]::xlSortNormal) | Out-Null
$ws.Sort.SortFields.Add($tieKey, [Microsoft.Office.Interop.Excel.XlSortOn]::xlSortOnValues, [Microsoft.Office.Interop.Excel.XlSortOrder]::xlAscending, [Microsoft.Office.Interop.Excel.XlSortDataOption]::xlSortNormal) | Out-Null
$ws.Sort.SetRange($dataRange)
$ws.Sort.Header = [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
$ws.Sort.Apply()

# --- Turn the data range into a proper Excel Table ---
$lo = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $dataRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$lo.Name = "Table1"
$lo.TableStyle = "TableStyleMedium2"

# --- View changes: zoom out and move the selection ---
$excel.ActiveWindow.Zoom = 160
$ws.Range("C14").Select()
